$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# S03/G01: Strategy service: metadata and parameter storage
# Mark the three tasks as implemented and add remarks describing the work done.

$ws.Range("G20").Value = "implemented"
$ws.Range("F20").Value = "Strategy CRUD API endpoints implemented."

$ws.Range("G21").Value = "implemented"
$ws.Range("F21").Value = "Strategy parameter CRUD endpoints wired to strategies."

$ws.Range("G22").Value = "implemented"
$ws.Range("F22").Value = "Tags/category/status/integration fields exposed in API models."
